$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.036.42'
$ws.Range("E2").Value = '  +1.92%  '

$ws.Range("D3").Value = '1.776.95'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.36'
$ws.Range("E5").Value = '  -0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3817'
$ws.Range("E7").Value = '  -2.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3419'
$ws.Range("E8").Value = '  -1.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.88'
$ws.Range("E9").Value = '  -2.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.141'
$ws.Range("E10").Value = '  -4.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07389'
$ws.Range("E11").Value = '  -1.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.30'
$ws.Range("E12").Value = '  +6.37%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  +0.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.383'
$ws.Range("E14").Value = '  -1.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.416'
$ws.Range("E15").Value = '  +3.71%  '

$ws.Range("D16").Value = '1.777.97'
$ws.Range("E16").Value = '  -1.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001074'
$ws.Range("E17").Value = '  -2.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06661'
$ws.Range("E18").Value = '  -0.50%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.26'
$ws.Range("E19").Value = '  -3.13%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.38%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.35'
$ws.Range("E21").Value = '  -1.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.411'
$ws.Range("E22").Value = '  -2.16%  '

$ws.Range("D23").Value = '28.078.13'
$ws.Range("E23").Value = '  +2.03%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.08'
$ws.Range("E24").Value = '  -2.84%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.380'
$ws.Range("E25").Value = '  -0.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.449'
$ws.Range("E26").Value = '  -0.94%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.72'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.412'
$ws.Range("E28").Value = '  -3.57%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.32'
$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("D30").Value = '1.980.40'
$ws.Range("E30").Value = '  -0.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.21'
$ws.Range("E31").Value = '  -0.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.025'
$ws.Range("E32").Value = '  -0.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.062'
$ws.Range("E33").Value = '  +0.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08880'
$ws.Range("E34").Value = '  +1.16%  '

$ws.Range("E35").Value = '  -2.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02396'
$ws.Range("E36").Value = '  -0.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6838'
$ws.Range("E37").Value = '  +0.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06407'
$ws.Range("E38").Value = '  -0.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.279'
$ws.Range("E39").Value = '  -3.15%  '

$ws.Range("E40").Value = '  -2.23%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.236'
$ws.Range("E41").Value = '  -1.34%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.498'
$ws.Range("E42").Value = '  -7.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.218'
$ws.Range("E43").Value = '  -1.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.27'
$ws.Range("E44").Value = '  -1.82%  '

$ws.Range("E45").Value = '  +0.33%  '

$ws.Range("E46").Value = '  -2.15%  '

$ws.Range("E47").Value = '  -0.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.82'
$ws.Range("E48").Value = '  +0.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.067'
$ws.Range("E49").Value = '  -3.30%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07499'
$ws.Range("E50").Value = '  +4.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.217'
$ws.Range("E51").Value = '  +5.01%  '
